$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.565.57'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '1.878.68'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  -1.23%  '
$ws.Range('D5').Value = '315.57'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  -1.34%  '
$ws.Range('D7').Value = '0.5110'
$ws.Range('E7').Value = '  -0.72%  '
$ws.Range('D8').Value = '0.3935'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('D9').Value = '0.08412'
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('D10').Value = '1.114'
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('D11').Value = '41.70'
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').Value = '6.285'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').Value = '1.882.03'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').Value = '20.51'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').Value = '7.288'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').Value = '1.006'
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').Value = '0.00001108'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '91.41'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').Value = '0.06722'
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('D20').Value = '17.75'
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = '1.006'
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').Value = '28.609.97'
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').Value = '11.14'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '2.251'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').Value = '2.097.10'
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('D27').Value = '161.04'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').Value = '2.381'
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('D30').Value = '126.82'
$ws.Range('E30').Value = '  +0.83%  '
$ws.Range('D31').Value = '0.1054'
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').Value = '3.612'
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').Value = '0.02464'
$ws.Range('E35').Value = '  +0.90%  '
$ws.Range('D36').Value = '0.06537'
$ws.Range('D37').Value = '0.2188'
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').Value = '8.963'
$ws.Range('E38').Value = '  -3.96%  '
$ws.Range('D39').Value = '1.261'
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('D40').Value = '1.200'
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.6488'
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D42').Value = '5.094'
$ws.Range('E42').Value = '  +2.45%  '
$ws.Range('D43').Value = '11.20'
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').Value = '1.006'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').Value = '0.6081'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = '13.09'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').Value = '3.703'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('D48').Value = '2.046'
$ws.Range('E48').Value = '  +2.04%  '
$ws.Range('D49').Value = '1.219'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '122.63'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('D51').Value = '1.187'
$ws.Range('E51').Value = '  -8.11%  '
